# Generate Report for Handoff
# Inserts a new data row (for file 09458abb-652e-48e8-8e5d-9581ec1c5232.md) above the
# existing row (4f014fc9-9bb1-48cf-b7b7-ec71c489cbe8.md) on all three worksheets
# (Overview, zh-cn, de-de), pushing the old row down to row 3, resizing the tables,
# and rewiring the hyperlinks so that:
#   - the "new" row occupies row 2 and links to the 09458abb file
#   - the "old" row is preserved at row 3 and links to the 4f014fc9 file

$wb = $excel.ActiveWorkbook

$commit = "b00192cad76c301c57e0009b9ca45eee0d1e8dcd"
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/"

$newFile = "09458abb-652e-48e8-8e5d-9581ec1c5232.md"
$oldFile = "4f014fc9-9bb1-48cf-b7b7-ec71c489cbe8.md"

function Set-HyperlinkStyle($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1.xml) — columns A:G
#   A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#   E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

# Row 2: new file
$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = "e2e\" + $newFile
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("D2").Value = "'"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 10:43:54"
$wsOverview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3: previous row data, shifted down
$wsOverview.Range("A3").Value = $oldFile
$wsOverview.Range("B3").Value = "e2e\" + $oldFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 10:43:38"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), ($repoBase + $newFile), [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $newFile)) | Out-Null
Set-HyperlinkStyle($wsOverview.Range("B2"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($repoBase + $oldFile), [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $oldFile)) | Out-Null
Set-HyperlinkStyle($wsOverview.Range("B3"))

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2.xml) and "de-de" (sheet3.xml) — columns A:P
#   A=Source File Name, B=File Extension, C=Status, D=Source Path, E=Priority,
#   F=Content Duplicate, G=Latest Handoff File, H=Latest Handoff Datetime,
#   I=Latest Target File, J=Latest Handback File, K=Latest Handback DateTime,
#   L=Reference Tokens, M=To be localized, N=Dependency From, O=Has metadata,
#   P=Error Detail
# ---------------------------------------------------------------------------
function Fill-LangSheet($ws, $newXlf, $newDate, $oldXlf, $oldDate) {
    $ws.Hyperlinks.Delete()

    # Row 2: new file
    $ws.Range("A2").Value = $newFile
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("D2").Value = "e2e"
    $ws.Range("E2").Value = "ht"
    $ws.Range("F2").Value = "'False"
    $ws.Range("G2").Value = $newXlf
    $ws.Range("H2").Value = $newDate
    $ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("I2").Value = "'"
    $ws.Range("J2").Value = "'"
    $ws.Range("K2").Value = "0001-01-01 00:00:00"
    $ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("L2").Value = "'"
    $ws.Range("M2").Value = "'True"
    $ws.Range("N2").Value = "'"
    $ws.Range("O2").Value = "'False"
    $ws.Range("P2").Value = "'"

    # Row 3: previous row data, shifted down
    $ws.Range("A3").Value = $oldFile
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = "e2e"
    $ws.Range("E3").Value = "ht"
    $ws.Range("F3").Value = "'False"
    $ws.Range("G3").Value = $oldXlf
    $ws.Range("H3").Value = $oldDate
    $ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("I3").Value = "'"
    $ws.Range("J3").Value = "'"
    $ws.Range("K3").Value = "0001-01-01 00:00:00"
    $ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("L3").Value = "'"
    $ws.Range("M3").Value = "'True"
    $ws.Range("N3").Value = "'"
    $ws.Range("O3").Value = "'False"
    $ws.Range("P3").Value = "'"

    $ws.Hyperlinks.Add($ws.Range("A2"), ($repoBase + $newFile), [System.Type]::Missing, [System.Type]::Missing, $newFile) | Out-Null
    Set-HyperlinkStyle($ws.Range("A2"))
    $ws.Hyperlinks.Add($ws.Range("A3"), ($repoBase + $oldFile), [System.Type]::Missing, [System.Type]::Missing, $oldFile) | Out-Null
    Set-HyperlinkStyle($ws.Range("A3"))

    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A1:P3"))
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Fill-LangSheet $wsZhCn `
    "09458abb-652e-48e8-8e5d-9581ec1c5232.31d4f97bc0e5eeeb11e0dc9f635ddbb294c944ef.zh-cn.xlf" "2016-08-18 10:43:48" `
    "4f014fc9-9bb1-48cf-b7b7-ec71c489cbe8.a577ee43e37bbd4ee1fed743fb049211baa1f20c.zh-cn.xlf" "2016-08-18 10:43:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
Fill-LangSheet $wsDeDe `
    "09458abb-652e-48e8-8e5d-9581ec1c5232.31d4f97bc0e5eeeb11e0dc9f635ddbb294c944ef.de-de.xlf" "2016-08-18 10:43:54" `
    "4f014fc9-9bb1-48cf-b7b7-ec71c489cbe8.a577ee43e37bbd4ee1fed743fb049211baa1f20c.de-de.xlf" "2016-08-18 10:43:38"

Write-Host "Done applying handback report update"
